# Auto-generated edit script: apply cell value updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.62%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '5'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.27'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.52%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '5'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.127'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.47%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '5'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05587'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.32%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '5'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.470'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.40%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '5'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8188'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.01%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '5'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8360'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.09%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '5'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1334'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.42%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '5'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06995'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.64%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '5'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02890'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.53%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '5'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.06%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '5'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001519'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.41%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '5'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006001'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.84%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '5'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006180'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.51%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '5'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.654'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '5'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.035'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.57%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '5'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '5'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.12%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '5'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03099'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.86%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '5'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.22%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '5'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.761'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.06%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '5'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04666'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.42%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '5'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.07%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '5'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001247'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.22%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '5'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.02%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '5'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-1.03%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '5'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '5'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '5'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '5'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '5'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '5'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '5'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '5'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '5'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '5'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '5'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '5'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '5'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03641'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.56%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '5'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1371'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.12%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '5'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002551'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.05%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '5'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003432'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-44.57%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '5'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008845'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.07%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '5'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005335'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.80%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '5'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.00%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '5'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '5'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10.22%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '5'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.00%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '5'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.00%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '5'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '5'
